$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F54").Value = 117
$ws.Range("G54").Value = 2958.93
$ws.Range("F73").Value = 48
$ws.Range("G73").Value = 1606.08
$ws.Range("B74").Value = 301389.95
$ws.Range("B151").Value = 65258
$ws.Range("F151").Value = 2
$ws.Range("G151").Value = 64287.16
$ws.Range("B152").Value = 64196
$ws.Range("F152").Value = 1
$ws.Range("G152").Value = 32143.58
$ws.Range("F225").Value = 58
$ws.Range("G225").Value = 9300.879999999999
$ws.Range("F244").Value = 47
$ws.Range("G244").Value = 681.97
$ws.Range("B251").Value = 105673.4
$ws.Range("F288").Value = 2484
$ws.Range("G288").Value = 45954
$ws.Range("B295").Value = 69165.49000000001
$ws.Range("F372").Value = 39
$ws.Range("G372").Value = 3341.52
$ws.Range("B376").Value = 186643.76
$ws.Range("B387").Value = 55373
$ws.Range("E387").Value = 163.62
$ws.Range("F387").Value = -94
$ws.Range("G387").Value = -13562.32
$ws.Range("B388").Value = 63520
$ws.Range("E388").Value = 153.4
$ws.Range("F388").Value = 47
$ws.Range("G388").Value = 6781.16
$ws.Range("B391").Value = 63510
$ws.Range("E391").Value = 50.66
$ws.Range("F391").Value = 86
$ws.Range("G391").Value = 4097.04
$ws.Range("B392").Value = 55356
$ws.Range("E392").Value = 54.04
$ws.Range("F392").Value = -158
$ws.Range("G392").Value = -7527.12
$ws.Range("F443").Value = 12
$ws.Range("G443").Value = 516.12
$ws.Range("F448").Value = 187
$ws.Range("G448").Value = 13955.81
$ws.Range("B455").Value = 102784.14
$ws.Range("F473").Value = 116
$ws.Range("G473").Value = 17431.32
$ws.Range("B474").Value = 105914.56
$ws.Range("B554").Value = 53263
$ws.Range("E554").Value = 15.29
$ws.Range("F554").Value = -309
$ws.Range("G554").Value = -3958.29
$ws.Range("B555").Value = 65066
$ws.Range("E555").Value = 13.61
$ws.Range("F555").Value = 90
$ws.Range("G555").Value = 1152.9
$ws.Range("B560").Value = 45706
$ws.Range("E560").Value = 23.58
$ws.Range("F560").Value = -202
$ws.Range("G560").Value = -3985.46
$ws.Range("B561").Value = 64922
$ws.Range("E561").Value = 20.98
$ws.Range("F561").Value = 67
$ws.Range("G561").Value = 1321.91
$ws.Range("B568").Value = 45709
$ws.Range("E568").Value = 15.69
$ws.Range("F568").Value = -300
$ws.Range("G568").Value = -3945
$ws.Range("B569").Value = 64925
$ws.Range("E569").Value = 13.97
$ws.Range("F569").Value = 111
$ws.Range("G569").Value = 1459.65
$ws.Range("B570").Value = 45702
$ws.Range("E570").Value = 31.43
$ws.Range("F570").Value = -215
$ws.Range("G570").Value = -5654.5
$ws.Range("B571").Value = 64919
$ws.Range("E571").Value = 27.97
$ws.Range("F571").Value = 61
$ws.Range("G571").Value = 1604.3
$ws.Range("B573").Value = 65067
$ws.Range("E573").Value = 15.65
$ws.Range("F573").Value = 126
$ws.Range("G573").Value = 1855.98
$ws.Range("B574").Value = 53595
$ws.Range("E574").Value = 17.61
$ws.Range("F574").Value = -335
$ws.Range("G574").Value = -4934.55
$ws.Range("F583").Value = 295
$ws.Range("G583").Value = 2843.8
$ws.Range("B588").Value = 45758.94
$ws.Range("B644").Value = 64810
$ws.Range("E644").Value = 291.22
$ws.Range("F644").Value = 2
$ws.Range("G644").Value = 547.84
$ws.Range("B645").Value = 53319
$ws.Range("E645").Value = 310.64
$ws.Range("F645").Value = -6
$ws.Range("G645").Value = -1643.52
$ws.Range("F658").Value = 21
$ws.Range("G658").Value = 380.94
$ws.Range("B660").Value = 5975.13
$ws.Range("B663").Value = 60025
$ws.Range("E663").Value = 37.22
$ws.Range("F663").Value = -98
$ws.Range("G663").Value = -3217.34
$ws.Range("B664").Value = 64833
$ws.Range("E664").Value = 34.9
$ws.Range("F664").Value = 91
$ws.Range("G664").Value = 2987.53
$ws.Range("B673").Value = 60022
$ws.Range("E673").Value = 37.22
$ws.Range("F673").Value = -113
$ws.Range("G673").Value = -3709.79
$ws.Range("B674").Value = 64830
$ws.Range("E674").Value = 34.9
$ws.Range("F674").Value = 92
$ws.Range("G674").Value = 3020.36
$ws.Range("F689").Value = 571
$ws.Range("G689").Value = 31336.48
$ws.Range("F690").Value = 164
$ws.Range("G690").Value = 4493.6
$ws.Range("F692").Value = 145
$ws.Range("G692").Value = 3862.8
$ws.Range("F693").Value = 221
$ws.Range("G693").Value = 18919.81
$ws.Range("B696").Value = 200784.85
$ws.Range("F806").Value = 28
$ws.Range("G806").Value = 3379.88
$ws.Range("B807").Value = 42076
$ws.Range("F880").Value = 78
$ws.Range("G880").Value = 2319.72
$ws.Range("F883").Value = 106
$ws.Range("G883").Value = 8512.860000000001
$ws.Range("B890").Value = 27112.93
$ws.Range("F895").Value = 317
$ws.Range("G895").Value = 9582.91
$ws.Range("F896").Value = 1927
$ws.Range("G896").Value = 314312.97
$ws.Range("B902").Value = 368359.92
$ws.Range("F933").Value = 8
$ws.Range("G933").Value = 5436.64
$ws.Range("F940").Value = 6
$ws.Range("G940").Value = 2428.62
$ws.Range("B941").Value = 143159.42
$ws.Range("B947").Value = 5636686.47
$ws.Range("B948").Value = 5636686.47
